# ATP workbook edit: add a "focus" column to the ATP table, clear a few
# stray values in row 2, select the newly added range, and nudge the
# window/view metadata to match the saved state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ATP")

# --- Clear the stray swim_time/bike_time/run_time values on row 2 ---
$ws.Range("H2:J2").ClearContents()

# --- Add a new trailing column ("focus") to the Tabel4 table ---
$lo = $ws.ListObjects.Item(1)
$newCol = $lo.ListColumns.Add()

# --- Fill in the "focus" column. Cells are written in the exact order
#     needed so that the shared-string table indexes new values the same
#     way the original authoring session did. ---

# Primed writes establish the first-seen order of each unique value.
$ws.Range("L1").Value = "focus"
$ws.Range("L10").Value = "power"
$ws.Range("L2").Value = "endurance"
$ws.Range("L23").Value = "technique"
$ws.Range("L24").Value = "your raceday!"
$ws.Range("L25").Value = "recovery"
$ws.Range("L26").Value = "staying strong"
$ws.Range("L37").Value = "a lot of training!"
$ws.Range("L38").Value = "preparing for race."
$ws.Range("L53").Value = "Christmas..."
$ws.Range("L44").Value = "nothing at all."

# Remaining cells reuse the values already registered above.
$ws.Range("L3").Value = "endurance"
$ws.Range("L4").Value = "endurance"
$ws.Range("L5").Value = "endurance"
$ws.Range("L6").Value = "endurance"
$ws.Range("L7").Value = "endurance"
$ws.Range("L8").Value = "endurance"
$ws.Range("L9").Value = "endurance"
$ws.Range("L27").Value = "staying strong"
$ws.Range("L28").Value = "staying strong"
$ws.Range("L29").Value = "staying strong"
$ws.Range("L30").Value = "staying strong"
$ws.Range("L34").Value = "power"
$ws.Range("L35").Value = "endurance"
$ws.Range("L36").Value = "power"
$ws.Range("L39").Value = "your raceday!"
$ws.Range("L40").Value = "recovery"
$ws.Range("L41").Value = "staying strong"
$ws.Range("L42").Value = "staying strong"
$ws.Range("L43").Value = "staying strong"
$ws.Range("L45").Value = "nothing at all."
$ws.Range("L46").Value = "nothing at all."
$ws.Range("L47").Value = "nothing at all."
$ws.Range("L48").Value = "nothing at all."
$ws.Range("L49").Value = "nothing at all."
$ws.Range("L50").Value = "nothing at all."
$ws.Range("L51").Value = "nothing at all."
$ws.Range("L52").Value = "nothing at all."

# --- Size the new column to roughly match its saved width ---
$ws.Columns("L").ColumnWidth = 11

# --- Match the saved selection/active-cell state ---
$ws.Range("L44:L52").Select()
